$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.896.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.621.95"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.502"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.07%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0615"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.25"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.848.33"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.622.88"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.60%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.523"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.874.96"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.56%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.01"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.24%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.23"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.88%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.55"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.06"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.36%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.51"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.73"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.70"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.12"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.57%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.54%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.88%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.09"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.08%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.18%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.115.08"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.844"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.68%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.515"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.56%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.81"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.24%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "RocketPoolETH"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.758.45"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.11%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.763"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.41%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.89%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.96%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0529"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.17"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.24%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.14%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.23%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.41"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.77%  "
